$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row relabel + new "Approach" column header ---
$ws.Range("E2").Value = "Category"
$ws.Range("F2").Value = "Patterns"
$ws.Range("I2").Value = "Approach"

# --- New row 6 data: LeetCode 191 "Number of 1 bits" entry ---
$ws.Range("C6").Value = 191
$ws.Range("D6").Value = "Number of 1 bits "
$ws.Range("E6").Value = "Math"
$ws.Range("F6").Value = "binary search"

# H6: related-problems text with a couple of hyperlink-styled (underlined) runs
$h6 = $ws.Range("H6")
$h6.Value = "190. Reverse Bits,  231. Power of Two , 338. Counting Bits"
$h6.WrapText = $true

$run1 = $h6.Characters(21, 18)
$run1.Font.Underline = $true
$run1.Font.Color = 8812614

$run2 = $h6.Characters(41, 18)
$run2.Font.Underline = $true
$run2.Font.Color = 8812614

# I6: approach explanation text
$i6 = $ws.Range("I6")
$i6.Value = "The idea is to remove the rightmost one from n's binary. Which can be achieved by Subtracting 1 from n and perform bitwise operation between n and n-1."
$i6.WrapText = $true

# --- Row height for the newly wrapped row ---
$ws.Rows.Item(6).RowHeight = 43.2

# --- Column widths for the newly used columns ---
$ws.Range("E1").ColumnWidth = 10.917
$ws.Range("I1").ColumnWidth = 57.25

# --- View tweaks (zoom + final selection) ---
$excel.ActiveWindow.Zoom = 71
$ws.Range("I11").Select()
